# Weekly rolling update of price data (rows 173-206 in Sheet1).
# Every row's D/J/K/L/M/P values move "up" by two rows (row[i] <- row[i-2]),
# two brand-new rows are appended at the bottom (copies of the former last
# week, rows 205/206), and a brand-new week's data is written into the
# first two data rows (173/174).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 173
$lastRow = 206

# Column indexes: D=4, J=10, K=11, L=12, M=13, P=16
$colD = 4
$colJ = 10
$colK = 11
$colL = 12
$colM = 13
$colP = 16

# 1) Snapshot the current (pre-edit) values for the columns that move,
#    for every row in the block, before any writes happen.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, $colD).Value2
        J = $ws.Cells.Item($r, $colJ).Value2
        K = $ws.Cells.Item($r, $colK).Value2
        L = $ws.Cells.Item($r, $colL).Value2
        M = $ws.Cells.Item($r, $colM).Value2
        P = $ws.Cells.Item($r, $colP).Value2
    }
}

# 2) Append two brand-new rows (207/208) that duplicate the former last
#    week's two rows (205/206) exactly, before those rows get overwritten.
$ws.Range("A205:R205").Copy($ws.Range("A207:R207"))
$ws.Range("A206:R206").Copy($ws.Range("A208:R208"))

# 3) Shift every row's D/J/K/L/M/P values up by two rows: new row[r] takes
#    the pre-edit values that used to live at row (r-2).
for ($r = $lastRow; $r -ge ($firstRow + 2); $r--) {
    $src = $snapshot[$r - 2]
    $ws.Cells.Item($r, $colD).Value2 = $src.D
    $ws.Cells.Item($r, $colJ).Value2 = $src.J
    $ws.Cells.Item($r, $colK).Value2 = $src.K
    $ws.Cells.Item($r, $colL).Value2 = $src.L
    $ws.Cells.Item($r, $colM).Value2 = $src.M
    $ws.Cells.Item($r, $colP).Value2 = $src.P
}

# 4) Write the brand-new week's data into the first two data rows.
$ws.Cells.Item(173, $colD).Value2 = 44543
$ws.Cells.Item(173, $colJ).Value2 = 1200
$ws.Cells.Item(173, $colK).Value2 = 300
$ws.Cells.Item(173, $colL).Value2 = 350
$ws.Cells.Item(173, $colM).Value2 = 325
$ws.Cells.Item(173, $colP).Value2 = 81

$ws.Cells.Item(174, $colD).Value2 = 44543
$ws.Cells.Item(174, $colJ).Value2 = 1200
$ws.Cells.Item(174, $colK).Value2 = 300
$ws.Cells.Item(174, $colL).Value2 = 350
$ws.Cells.Item(174, $colM).Value2 = 325
$ws.Cells.Item(174, $colP).Value2 = 65

Write-Output "Done"
